# edit.ps1 — apply cryptos list update (commit: "Updated cryptos list on Fri Mar 24 21:22:22 UTC 2023 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All Price (column D) cells are formatted as Text first so that values such as
# "27.543.80" (multi-dot) or "1.400" (trailing zero) are preserved exactly as
# literal strings instead of being reinterpreted as numbers by Excel.
$priceCells = @(
    "D2",
    "D3",
    "D5",
    "D7",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D15",
    "D16",
    "D17",
    "D18",
    "D19",
    "D21",
    "D22",
    "D23",
    "D24",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Rows whose Coin/Link stay the same: update Price (D) and/or Volume(1h) (E) ---
$ws.Range("D2").Value = "27.543.80"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").Value = "1.750.90"
$ws.Range("E3").Value = "  -3.40%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "322.14"
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "0.4239"
$ws.Range("E7").Value = "  -3.69%  "
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D11").Value = "1.092"
$ws.Range("E11").Value = "  -2.72%  "
$ws.Range("D12").Value = "1.006"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "20.57"
$ws.Range("E13").Value = "  -6.45%  "
$ws.Range("D14").Value = "6.016"
$ws.Range("E14").Value = "  -3.66%  "
$ws.Range("D15").Value = "7.233"
$ws.Range("E15").Value = "  -4.06%  "
$ws.Range("D16").Value = "1.780.06"
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").Value = "90.62"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").Value = "0.00001069"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "0.06365"
$ws.Range("E19").Value = "  -3.98%  "
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").Value = "16.96"
$ws.Range("E21").Value = "  -3.09%  "
$ws.Range("D22").Value = "5.880"
$ws.Range("E22").Value = "  -5.10%  "
$ws.Range("D23").Value = "27.610.80"
$ws.Range("E23").Value = "  -2.17%  "
$ws.Range("D24").Value = "11.13"
$ws.Range("E24").Value = "  -4.53%  "
$ws.Range("E25").Value = "  +3.42%  "
$ws.Range("D26").Value = "161.13"
$ws.Range("E26").Value = "  +3.62%  "
$ws.Range("D27").Value = "20.17"
$ws.Range("E27").Value = "  -2.53%  "
$ws.Range("D28").Value = "1.971.39"
$ws.Range("E28").Value = "  -2.55%  "
$ws.Range("D29").Value = "2.124"
$ws.Range("E29").Value = "  -8.01%  "
$ws.Range("D30").Value = "124.50"
$ws.Range("E30").Value = "  -2.83%  "
$ws.Range("D31").Value = "1.101"
$ws.Range("E31").Value = "  -8.23%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").Value = "5.538"
$ws.Range("E33").Value = "  -5.33%  "
$ws.Range("D34").Value = "0.08864"
$ws.Range("E34").Value = "  -3.74%  "
$ws.Range("D35").Value = "12.19"
$ws.Range("E35").Value = "  -6.44%  "
$ws.Range("D36").Value = "0.02284"
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("D37").Value = "0.2096"
$ws.Range("E37").Value = "  -3.27%  "
$ws.Range("D38").Value = "0.06001"
$ws.Range("E38").Value = "  -3.34%  "
$ws.Range("D39").Value = "0.6326"
$ws.Range("E39").Value = "  -3.38%  "
$ws.Range("D40").Value = "4.931"
$ws.Range("E40").Value = "  -4.08%  "
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("D45").Value = "13.32"
$ws.Range("E45").Value = "  -4.33%  "
$ws.Range("D46").Value = "0.5856"
$ws.Range("E46").Value = "  -3.36%  "
$ws.Range("D47").Value = "3.692"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("D48").Value = "123.19"
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("D49").Value = "1.980"
$ws.Range("E49").Value = "  -2.44%  "
$ws.Range("D50").Value = "1.162"
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("E51").Value = "  -2.42%  "

# --- Rows 9/10 and 43/44: the two coins swapped rank/position, so Coin, Link, Price and Volume all change ---
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.07493"
$ws.Range("E9").Value = "  -2.44%  "

$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "42.44"
$ws.Range("E10").Value = "  -5.09%  "

$ws.Range("B43").Value = "WEMIXTOKEN"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "1.400"
$ws.Range("E43").Value = "  +1.01%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "7.844"
$ws.Range("E44").Value = "  -3.59%  "

